{"js": "// PlanIt Questions and Answers \u2014 \"added few more supporting points\"\n//\n// This applies the five substantive content changes from the commit:\n//   1. Login functionality bullet: trailing space added after \"...data validations etc.\"\n//   2. Contact page bullet: \"back button functionality, \" inserted before the\n//      trailing \"etc\"\n//   3. Cart Page bullet: \"back button after empty cart, \" inserted after\n//      \"Empty \"\n//   4. Q2 \"optimized code\" answer: \" such as give appropriate time for\n//      waits, different conditions etc.\" appended (replacing the final period)\n//   5. Framework answer: a new closing sentence about \"testing xml\" appended\n\nconst body = context.document.body;\n\n// 1. Login functionality: \"...data validations etc.\" -> \"...data validations etc. \"\nconst loginHits = body.search(\"data validations etc.\", { matchCase: true, matchWholeWord: false });\nloginHits.load(\"items\");\nawait context.sync();\nif (loginHits.items.length !== 1) {\n  throw new Error(\"Expected 1 match for login functionality text, got \" + loginHits.items.length);\n}\nloginHits.items[0].insertText(\"data validations etc. \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Contact page: \"...Cross Site Request Forgery tests, etc\" ->\n//    \"...Cross Site Request Forgery tests, back button functionality, etc\"\nconst contactHits = body.search(\"Cross Site Request Forgery tests, etc\", { matchCase: true });\ncontactHits.load(\"items\");\nawait context.sync();\nif (contactHits.items.length !== 1) {\n  throw new Error(\"Expected 1 match for contact page text, got \" + contactHits.items.length);\n}\ncontactHits.items[0].insertText(\n  \"Cross Site Request Forgery tests, back button functionality, etc\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 3. Cart Page: \"Empty cart, remove item\" -> \"Empty cart, back button after empty cart, remove item\"\nconst cartHits = body.search(\"Empty cart, remove item\", { matchCase: true });\ncartHits.load(\"items\");\nawait context.sync();\nif (cartHits.items.length !== 1) {\n  throw new Error(\"Expected 1 match for cart page text, got \" + cartHits.items.length);\n}\ncartHits.items[0].insertText(\n  \"Empty cart, back button after empty cart, remove item\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 4. Q2 optimized-code answer: \"...will be the key.\" ->\n//    \"...will be the key such as give appropriate time for waits, different conditions etc.\"\nconst keyHits = body.search(\"Apart from this writing optimized code for execution will be the key.\", { matchCase: true });\nkeyHits.load(\"items\");\nawait context.sync();\nif (keyHits.items.length !== 1) {\n  throw new Error(\"Expected 1 match for optimized code text, got \" + keyHits.items.length);\n}\nkeyHits.items[0].insertText(\n  \"Apart from this writing optimized code for execution will be the key such as give appropriate time for waits, different conditions etc.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 5. Framework answer: append new closing sentence right after\n//    \"...(this I have done in current company)\"\nconst frameworkHits = body.search(\"(this I have done in current company)\", { matchCase: true });\nframeworkHits.load(\"items\");\nawait context.sync();\nif (frameworkHits.items.length !== 1) {\n  throw new Error(\"Expected 1 match for framework answer text, got \" + frameworkHits.items.length);\n}\nframeworkHits.items[0].insertText(\n  \"(this I have done in current company). Using testing xml also we can execute the tests in parallel which will reduce the execution time.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# PlanIt Questions and Answers -- \"added few more supporting points\"\n#\n# Applies the five substantive content changes from the commit using\n# Find/Replace (Content.Find.Execute) against the whole document:\n#   1. Login functionality bullet: trailing space added after\n#      \"...data validations etc.\"\n#   2. Contact page bullet: \"back button functionality, \" inserted before\n#      the trailing \"etc\"\n#   3. Cart Page bullet: \"back button after empty cart, \" inserted after\n#      \"Empty \"\n#   4. Q2 \"optimized code\" answer: \" such as give appropriate time for\n#      waits, different conditions etc.\" appended (replacing the final period)\n#   5. Framework answer: a new closing sentence about \"testing xml\" appended\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\n# 1. Login functionality: \"...data validations etc.\" -> \"...data validations etc. \"\nReplace-Text \"data validations etc.\" \"data validations etc. \"\n\n# 2. Contact page: \"...Cross Site Request Forgery tests, etc\" ->\n#    \"...Cross Site Request Forgery tests, back button functionality, etc\"\nReplace-Text \"Cross Site Request Forgery tests, etc\" \"Cross Site Request Forgery tests, back button functionality, etc\"\n\n# 3. Cart Page: \"Empty cart, remove item\" -> \"Empty cart, back button after empty cart, remove item\"\nReplace-Text \"Empty cart, remove item\" \"Empty cart, back button after empty cart, remove item\"\n\n# 4. Q2 optimized-code answer: \"...will be the key.\" ->\n#    \"...will be the key such as give appropriate time for waits, different conditions etc.\"\nReplace-Text \"Apart from this writing optimized code for execution will be the key.\" \"Apart from this writing optimized code for execution will be the key such as give appropriate time for waits, different conditions etc.\"\n\n# 5. Framework answer: append new closing sentence right after\n#    \"...(this I have done in current company)\"\nReplace-Text \"(this I have done in current company)\" \"(this I have done in current company). Using testing xml also we can execute the tests in parallel which will reduce the execution time.\"\n"}
